$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq $name) { return $shp }
    }
    return $null
}

# --- "テキスト ボックス 46": "：読める" -> "：" + "見える" -----------------
$shp1 = Get-ShapeByName $s "テキスト ボックス 46"
$tr1 = $shp1.TextFrame.TextRange
$sub1 = $tr1.Characters(2, $tr1.Length - 1)
$sub1.Text = "見える"

# --- "テキスト ボックス 48": "：読めない" -> "：" + "見えない" -------------
$shp2 = Get-ShapeByName $s "テキスト ボックス 48"
$tr2 = $shp2.TextFrame.TextRange
$sub2 = $tr2.Characters(2, $tr2.Length - 1)
$sub2.Text = "見えない"

# --- "テキスト ボックス 55": 2nd paragraph "...みたい" -> "...見たい" -----
$shp3 = Get-ShapeByName $s "テキスト ボックス 55"
$tr3 = $shp3.TextFrame.TextRange
$para3 = $tr3.Paragraphs(2, 1)
$full3 = $para3.Characters(1, $para3.Length)
$full3.Text = "ホラーでない映画が見たい"
